$wb = $excel.ActiveWorkbook

# --- TestData sheet: duplicate row 2 (TC1 / Mercury / mercury) into rows 4, 5, 6 ---
$wsTestData = $wb.Worksheets.Item("TestData")
foreach ($r in 4..6) {
    $wsTestData.Range("A2:C2").Copy() | Out-Null
    $wsTestData.Range("A$r`:C$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $wsTestData.Range("A2:C2").Copy() | Out-Null
    $wsTestData.Range("A$r`:C$r").PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
$excel.CutCopyMode = $false

# --- Update the remembered cursor/selection position on each sheet ---
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsTestCases.Range("D12").Select() | Out-Null

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("E17").Select() | Out-Null

$wsTestData.Activate() | Out-Null
$wsTestData.Range("G13").Select() | Out-Null
